$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("H18").Value = 988.8889
$ws.Range("N18").Value = -1568
$ws.Range("I31").Value = 1305.875
$ws.Range("M31").Value = -3687.625
$ws.Range("H31").Value = 1305.875
$ws.Range("K31").Value = 3917.625
$ws.Range("J40").Value = 7426.857
$ws.Range("L40").Value = 7426.857
$ws.Range("H40").Value = 7123.5
$ws.Range("N40").Value = -7776.857
$ws.Range("I69").Value = 6657.4
$ws.Range("M69").Value = -19098.2
$ws.Range("H69").Value = 9759.571
$ws.Range("K69").Value = 19972.2
$ws.Range("I72").Value = 6657.4
$ws.Range("M72").Value = -55548.6
$ws.Range("H72").Value = 9759.571
$ws.Range("K72").Value = 59916.6
$ws.Range("I74").Value = 20415048
$ws.Range("M74").Value = -20414112
$ws.Range("H74").Value = 15881148
$ws.Range("K74").Value = 20415048
$ws.Range("I77").Value = 20415048
$ws.Range("M77").Value = -102070560
$ws.Range("H77").Value = 15881148
$ws.Range("K77").Value = 102075240
$ws.Range("I103").Value = 733.3333
$ws.Range("M103").Value = -1613.9999
$ws.Range("H103").Value = 733.3333
$ws.Range("K103").Value = 2199.9999
$ws.Range("I111").Value = 809.7778
$ws.Range("M111").Value = 637.6666
$ws.Range("H111").Value = 2253.4546
$ws.Range("K111").Value = 2429.3334
$ws.Range("J132").Value = 13979.8
$ws.Range("L132").Value = 41939.39999999999
$ws.Range("H132").Value = 493163.78
$ws.Range("N132").Value = -46999.39999999999
$ws.Range("I137").Value = 8654.286
$ws.Range("M137").Value = -23412.858
$ws.Range("H137").Value = 7332.615
$ws.Range("K137").Value = 25962.858
$ws.Range("I138").Value = 1846.6522
$ws.Range("M138").Value = -399.9565999999995
$ws.Range("H138").Value = 3382.9807
$ws.Range("K138").Value = 5539.9566

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 3784.585
$ws.Range("M32").Value = -3497.585
$ws.Range("J32").Value = 55569556
$ws.Range("L32").Value = 55569556
$ws.Range("H32").Value = 2024358.1
$ws.Range("K32").Value = 3784.585
$ws.Range("N32").Value = -55570130
$ws.Range("I61").Value = 8884.6
$ws.Range("M61").Value = -8672.6
$ws.Range("H61").Value = 5436.14
$ws.Range("K61").Value = 8884.6
$ws.Range("I74").Value = 9392.857
$ws.Range("M74").Value = -8518.857
$ws.Range("J74").Value = 4324.6665
$ws.Range("L74").Value = 4324.6665
$ws.Range("H74").Value = 7872.4
$ws.Range("K74").Value = 9392.857
$ws.Range("N74").Value = -6072.6665
$ws.Range("I77").Value = 9392.857
$ws.Range("M77").Value = -42596.285
$ws.Range("J77").Value = 4324.6665
$ws.Range("L77").Value = 21623.3325
$ws.Range("H77").Value = 7872.4
$ws.Range("K77").Value = 46964.285
$ws.Range("N77").Value = -30359.3325
$ws.Range("I88").Value = 5000.5
$ws.Range("M88").Value = -4594.5
$ws.Range("J88").Value = 6479
$ws.Range("L88").Value = 6479
$ws.Range("H88").Value = 5739.75
$ws.Range("K88").Value = 5000.5
$ws.Range("N88").Value = -7291
$ws.Range("I91").Value = 5000.5
$ws.Range("M91").Value = -3596.5
$ws.Range("J91").Value = 6479
$ws.Range("L91").Value = 6479
$ws.Range("H91").Value = 5739.75
$ws.Range("K91").Value = 5000.5
$ws.Range("N91").Value = -9287
$ws.Range("I110").Value = 1873.4706
$ws.Range("M110").Value = 171.5293999999999
$ws.Range("H110").Value = 4396.6875
$ws.Range("K110").Value = 1873.4706
$ws.Range("I136").Value = 8884.6
$ws.Range("M136").Value = -24103.8
$ws.Range("H136").Value = 5436.14
$ws.Range("K136").Value = 26653.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I80").Value = 3333
$ws.Range("M80").Value = -2335
$ws.Range("J80").Value = 37064560
$ws.Range("L80").Value = 37064560
$ws.Range("H80").Value = 27799252
$ws.Range("K80").Value = 3333
$ws.Range("N80").Value = -37066556
$ws.Range("I83").Value = 3333
$ws.Range("M83").Value = -11673
$ws.Range("J83").Value = 37064560
$ws.Range("L83").Value = 185322800
$ws.Range("H83").Value = 27799252
$ws.Range("K83").Value = 16665
$ws.Range("N83").Value = -185332784
$ws.Range("I86").Value = 1267.2084
$ws.Range("M86").Value = -144.2084
$ws.Range("H86").Value = 2797
$ws.Range("K86").Value = 1267.2084
$ws.Range("I89").Value = 1267.2084
$ws.Range("M89").Value = -720.0419999999995
$ws.Range("H89").Value = 2797
$ws.Range("K89").Value = 6336.041999999999
$ws.Range("I128").Value = 4439.5
$ws.Range("M128").Value = -10828.5
$ws.Range("H128").Value = 4439.5
$ws.Range("K128").Value = 13318.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 38477308
$ws.Range("M31").Value = -38477013
$ws.Range("J31").Value = 6291.143
$ws.Range("L31").Value = 6291.143
$ws.Range("H31").Value = 14715797
$ws.Range("K31").Value = 38477308
$ws.Range("N31").Value = -6881.143
$ws.Range("I34").Value = 38477308
$ws.Range("M34").Value = -38477106
$ws.Range("J34").Value = 6291.143
$ws.Range("L34").Value = 6291.143
$ws.Range("H34").Value = 14715797
$ws.Range("K34").Value = 38477308
$ws.Range("N34").Value = -6695.143
$ws.Range("I62").Value = 8721.223
$ws.Range("M62").Value = -8097.223
$ws.Range("H62").Value = 8024.091
$ws.Range("K62").Value = 8721.223
$ws.Range("I65").Value = 8721.223
$ws.Range("M65").Value = -40486.115
$ws.Range("H65").Value = 8024.091
$ws.Range("K65").Value = 43606.115
$ws.Range("I107").Value = 256.875
$ws.Range("M107").Value = 1663.125
$ws.Range("J107").Value = 1216.4
$ws.Range("L107").Value = 1216.4
$ws.Range("H107").Value = 625.9231
$ws.Range("K107").Value = 256.875
$ws.Range("N107").Value = -5056.4
$ws.Range("I122").Value = 1898.8
$ws.Range("M122").Value = -3246.4
$ws.Range("J122").Value = 20333
$ws.Range("L122").Value = 60999
$ws.Range("H122").Value = 6152.846
$ws.Range("K122").Value = 5696.4
$ws.Range("N122").Value = -65899
$ws.Range("I132").Value = 3896.1853
$ws.Range("M132").Value = -9158.555899999999
$ws.Range("H132").Value = 7656.467
$ws.Range("K132").Value = 11688.5559
$ws.Range("J141").Value = 416823.16
$ws.Range("L141").Value = 416823.16
$ws.Range("H141").Value = 365134.16
$ws.Range("N141").Value = -427183.16

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("K133").Value = 0

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I17").Value = 499
$ws.Range("M17").Value = -331
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("H17").Value = 499
$ws.Range("K17").Value = 499
$ws.Range("N17").ClearContents()
$ws.Range("I80").Value = 2444.75
$ws.Range("M80").Value = -1446.75
$ws.Range("J80").Value = 10338.223
$ws.Range("L80").Value = 10338.223
$ws.Range("H80").Value = 7909.4614
$ws.Range("K80").Value = 2444.75
$ws.Range("N80").Value = -12334.223
$ws.Range("I83").Value = 2444.75
$ws.Range("M83").Value = -7231.75
$ws.Range("J83").Value = 10338.223
$ws.Range("L83").Value = 51691.115
$ws.Range("H83").Value = 7909.4614
$ws.Range("K83").Value = 12223.75
$ws.Range("N83").Value = -61675.115
$ws.Range("I97").Value = 3128.5715
$ws.Range("M97").Value = -2632.5715
$ws.Range("H97").Value = 3269.9
$ws.Range("K97").Value = 3128.5715
$ws.Range("I113").Value = 2236.3333
$ws.Range("M113").Value = -66.33329999999978
$ws.Range("J113").Value = 12748.375
$ws.Range("L113").Value = 12748.375
$ws.Range("H113").Value = 8243.214
$ws.Range("K113").Value = 2236.3333
$ws.Range("N113").Value = -17088.375
$ws.Range("I126").Value = 83336056
$ws.Range("M126").Value = -250005698
$ws.Range("J126").Value = 24832.5
$ws.Range("L126").Value = 74497.5
$ws.Range("H126").Value = 41680444
$ws.Range("K126").Value = 250008168
$ws.Range("N126").Value = -79437.5
$ws.Range("I132").Value = 7435.15
$ws.Range("M132").Value = -19775.45
$ws.Range("J132").Value = 5362.125
$ws.Range("L132").Value = 16086.375
$ws.Range("H132").Value = 6842.857
$ws.Range("K132").Value = 22305.45
$ws.Range("N132").Value = -21146.375
$ws.Range("J134").Value = 66203.57000000001
$ws.Range("L134").Value = 198610.71
$ws.Range("H134").Value = 66203.57000000001
$ws.Range("N134").Value = -203680.71
$ws.Range("J136").Value = 22203
$ws.Range("L136").Value = 66609
$ws.Range("H136").Value = 22203
$ws.Range("N136").Value = -71709

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I55").Value = 5393
$ws.Range("M55").Value = -5220
$ws.Range("H55").Value = 5726.4116
$ws.Range("K55").Value = 5393
$ws.Range("I122").Value = 3994398.2
$ws.Range("M122").Value = -11980744.6
$ws.Range("J122").Value = 2934
$ws.Range("L122").Value = 8802
$ws.Range("H122").Value = 3073291.2
$ws.Range("K122").Value = 11983194.6
$ws.Range("N122").Value = -13702
$ws.Range("I132").Value = 47623050
$ws.Range("M132").Value = -142866620
$ws.Range("H132").Value = 35719172
$ws.Range("K132").Value = 142869150
$ws.Range("I136").Value = 13898139
$ws.Range("M136").Value = -41691867
$ws.Range("H136").Value = 12359457
$ws.Range("K136").Value = 41694417

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I136").Value = 55580296
$ws.Range("M136").Value = -166738338
$ws.Range("H136").Value = 27797910
$ws.Range("K136").Value = 166740888
$ws.Range("J140").Value = 83122.664
$ws.Range("L140").Value = 83122.664
$ws.Range("H140").Value = 83122.664
$ws.Range("N140").Value = -93482.664
